$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1980830670926517
$ws.Range("C2").Value = 0.5079872204472844
$ws.Range("J2").Value = 0.02236421725239617
$ws.Range("P2").Value = 0.1533546325878594
$ws.Range("S2").Value = 0.1182108626198083
$ws.Range("B3").Value = 0.025
$ws.Range("C3").Value = 0.0125
$ws.Range("J3").Value = 0.06875000000000001
$ws.Range("P3").Value = 0.70625
$ws.Range("S3").Value = 0.1875
$ws.Range("J4").Value = 0.07407407407407407
$ws.Range("P4").Value = 0.8148148148148148
$ws.Range("S4").Value = 0.1111111111111111
$ws.Range("B6").Value = 0.07650273224043716
$ws.Range("D6").Value = 0.00546448087431694
$ws.Range("F6").Value = 0.0273224043715847
$ws.Range("J6").Value = 0.2841530054644809
$ws.Range("O6").Value = 0.01092896174863388
$ws.Range("Q6").Value = 0.2076502732240437
$ws.Range("R6").Value = 0.06010928961748634
$ws.Range("S6").Value = 0.3278688524590164
$ws.Range("B7").Value = 0.1071428571428571
$ws.Range("D7").Value = 0.008928571428571428
$ws.Range("F7").Value = 0.04464285714285714
$ws.Range("J7").Value = 0.1651785714285714
$ws.Range("O7").Value = 0.02232142857142857
$ws.Range("Q7").Value = 0.1785714285714286
$ws.Range("R7").Value = 0.07589285714285714
$ws.Range("S7").Value = 0.3973214285714285
$ws.Range("B8").Value = 0.12
$ws.Range("D8").Value = 0.01538461538461539
$ws.Range("F8").Value = 0.05846153846153846
$ws.Range("J8").Value = 0.1538461538461539
$ws.Range("O8").Value = 0.03384615384615385
$ws.Range("Q8").Value = 0.1292307692307692
$ws.Range("R8").Value = 0.05230769230769231
$ws.Range("S8").Value = 0.4369230769230769
$ws.Range("B9").Value = 0.1428571428571428
$ws.Range("D9").Value = 0.006211180124223602
$ws.Range("F9").Value = 0.05590062111801242
$ws.Range("J9").Value = 0.1428571428571428
$ws.Range("O9").Value = 0.0124223602484472
$ws.Range("Q9").Value = 0.1490683229813665
$ws.Range("R9").Value = 0.06832298136645963
$ws.Range("S9").Value = 0.422360248447205
$ws.Range("B10").Value = 0.1206434316353887
$ws.Range("D10").Value = 0.01966041108132261
$ws.Range("F10").Value = 0.064343163538874
$ws.Range("J10").Value = 0.1492403932082216
$ws.Range("O10").Value = 0.03753351206434316
$ws.Range("Q10").Value = 0.2082216264521894
$ws.Range("R10").Value = 0.05898123324396783
$ws.Range("S10").Value = 0.3413762287756926
$ws.Range("G11").Value = 0.1613832853025937
$ws.Range("J11").Value = 0.0893371757925072
$ws.Range("K11").Value = 0.2046109510086455
$ws.Range("L11").Value = 0.5360230547550432
$ws.Range("S11").Value = 0.008645533141210375
$ws.Range("G12").Value = 0.7397959183673469
$ws.Range("J12").Value = 0.1836734693877551
$ws.Range("K12").Value = 0.01530612244897959
$ws.Range("L12").Value = 0.03571428571428571
$ws.Range("S12").Value = 0.02551020408163265
$ws.Range("F13").Value = 0.02631578947368421
$ws.Range("G13").Value = 0.7105263157894737
$ws.Range("J13").Value = 0.2368421052631579
$ws.Range("S13").Value = 0.02631578947368421
$ws.Range("J14").Value = 1
$ws.Range("F15").Value = 0.03045685279187817
$ws.Range("H15").Value = 0.1725888324873096
$ws.Range("I15").Value = 0.07106598984771574
$ws.Range("J15").Value = 0.3147208121827411
$ws.Range("K15").Value = 0.09644670050761421
$ws.Range("M15").Value = 0.01522842639593909
$ws.Range("O15").Value = 0.05583756345177665
$ws.Range("S15").Value = 0.2436548223350254
$ws.Range("F16").Value = 0.01129943502824859
$ws.Range("H16").Value = 0.1299435028248588
$ws.Range("I16").Value = 0.07344632768361582
$ws.Range("J16").Value = 0.480225988700565
$ws.Range("K16").Value = 0.1073446327683616
$ws.Range("M16").Value = 0.02824858757062147
$ws.Range("O16").Value = 0.03389830508474576
$ws.Range("S16").Value = 0.1355932203389831
$ws.Range("F17").Value = 0.01098901098901099
$ws.Range("H17").Value = 0.1428571428571428
$ws.Range("I17").Value = 0.09065934065934066
$ws.Range("J17").Value = 0.4313186813186813
$ws.Range("K17").Value = 0.1456043956043956
$ws.Range("M17").Value = 0.02747252747252747
$ws.Range("O17").Value = 0.03846153846153846
$ws.Range("S17").Value = 0.1126373626373626
$ws.Range("F18").Value = 0.01680672268907563
$ws.Range("H18").Value = 0.1176470588235294
$ws.Range("I18").Value = 0.1428571428571428
$ws.Range("J18").Value = 0.4285714285714285
$ws.Range("K18").Value = 0.1596638655462185
$ws.Range("M18").Value = 0.01680672268907563
$ws.Range("O18").Value = 0.008403361344537815
$ws.Range("S18").Value = 0.1092436974789916
$ws.Range("F19").Value = 0.03247373447946514
$ws.Range("H19").Value = 0.2015281757402101
$ws.Range("I19").Value = 0.07927411652340019
$ws.Range("J19").Value = 0.3371537726838587
$ws.Range("K19").Value = 0.1442215854823305
$ws.Range("M19").Value = 0.0171919770773639
$ws.Range("N19").Value = 0.0009551098376313276
$ws.Range("O19").Value = 0.07927411652340019
$ws.Range("S19").Value = 0.10792741165234
